# Developed new automation script (DRAIAM102)
#
# Adds a new test-case row (row 31) to the "Test Cases" sheet, following
# the same TCID / Jira id / Description / Runmode / Results layout used
# by every other row in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 carries the same per-column formatting (bordered cells, wrapped
# Description column, distinct Runmode cell format) that the new row
# needs, so clone its formats onto the new row first ...
$ws.Range("A28:E28").Copy()
$ws.Range("A31:E31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ... then fill in the new test case's data. Values are written in the
# same column order as the original authoring session (TCID, Description,
# Jira id, Runmode) so shared strings come out in that order.
$ws.Cells.Item(31, 1).Value = "DRAIAM102"
$ws.Cells.Item(31, 3).Value = "Verify that 'EndNote' should be moved within the white area and should be above 'Forgot Password' text and center aligned`n|| Verify that Clarivate Analytics logo should be Placed below the marketing area (centered)."
$ws.Cells.Item(31, 2).Value = "OPQA-5136 || OPQA-5137"
$ws.Cells.Item(31, 4).Value = "Y"

$ws.Rows.Item(31).RowHeight = 30

$ws.Range("B33").Select() | Out-Null
